$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 171.76923
$ws.Range("I5").Value = 79.5
$ws.Range("J5").Value = 319.4
$ws.Range("K5").Value = 79.5
$ws.Range("L5").Value = 319.4
$ws.Range("M5").Value = 35.5
$ws.Range("N5").Value = -549.4

$ws.Range("H33").Value = 226.17647
$ws.Range("I33").Value = 226.17647
$ws.Range("K33").Value = 226.17647
$ws.Range("M33").Value = 2.823530000000005

$ws.Range("H40").Value = 2923.7693
$ws.Range("I40").Value = 2143.5715
$ws.Range("K40").Value = 2143.5715
$ws.Range("M40").Value = -1968.5715

$ws.Range("H92").Value = 725.3333
$ws.Range("I92").Value = 336.81818
$ws.Range("K92").Value = 336.81818
$ws.Range("M92").Value = 911.18182

$ws.Range("H132").Value = 3667.075
$ws.Range("I132").Value = 3179.4849
$ws.Range("K132").Value = 9538.4547
$ws.Range("M132").Value = -7008.4547

$ws.Range("H137").Value = 7128675.5
$ws.Range("I137").Value = 938395.8
$ws.Range("K137").Value = 2815187.4
$ws.Range("M137").Value = -2812637.4

$ws.Range("H138").Value = 6110.5625
$ws.Range("J138").Value = 8714.950000000001
$ws.Range("L138").Value = 26144.85
$ws.Range("N138").Value = -36424.85000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 14823.5
$ws.Range("J46").Value = 14791.2
$ws.Range("L46").Value = 14791.2
$ws.Range("N46").Value = -15429.2

$ws.Range("H61").Value = 4823.8413
$ws.Range("I61").Value = 4740.1455
$ws.Range("J61").Value = 5399.25
$ws.Range("K61").Value = 4740.1455
$ws.Range("L61").Value = 5399.25
$ws.Range("M61").Value = -4528.1455
$ws.Range("N61").Value = -5823.25

$ws.Range("H136").Value = 4823.8413
$ws.Range("I136").Value = 4740.1455
$ws.Range("J136").Value = 5399.25
$ws.Range("K136").Value = 14220.4365
$ws.Range("L136").Value = 16197.75
$ws.Range("M136").Value = -11670.4365
$ws.Range("N136").Value = -21297.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 47026.4
$ws.Range("I105").Value = 54650.617
$ws.Range("K105").Value = 54650.617
$ws.Range("M105").Value = -52903.617

$ws.Range("H134").Value = 5535.641
$ws.Range("I134").Value = 5704.5586
$ws.Range("J134").Value = 4387
$ws.Range("K134").Value = 17113.6758
$ws.Range("L134").Value = 13161
$ws.Range("M134").Value = -14578.6758
$ws.Range("N134").Value = -18231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2285.6333
$ws.Range("I31").Value = 1464.8667
$ws.Range("K31").Value = 1464.8667
$ws.Range("M31").Value = -1169.8667

$ws.Range("H34").Value = 2285.6333
$ws.Range("I34").Value = 1464.8667
$ws.Range("K34").Value = 1464.8667
$ws.Range("M34").Value = -1262.8667

$ws.Range("H107").Value = 26484.385
$ws.Range("I107").Value = 37033
$ws.Range("J107").Value = 2750
$ws.Range("K107").Value = 37033
$ws.Range("L107").Value = 2750
$ws.Range("M107").Value = -35113
$ws.Range("N107").Value = -6590

$ws.Range("H122").Value = 2230.8
$ws.Range("I122").Value = 1788.75
$ws.Range("K122").Value = 5366.25
$ws.Range("M122").Value = -2916.25

$ws.Range("H132").Value = 20867430
$ws.Range("I132").Value = 27789072
$ws.Range("K132").Value = 83367216
$ws.Range("M132").Value = -83364686

$ws.Range("H134").Value = 1529517.9
$ws.Range("I134").Value = 2410018.5
$ws.Range("K134").Value = 7230055.5
$ws.Range("M134").Value = -7227520.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1534183.4
$ws.Range("I34").Value = 3126471
$ws.Range("J34").Value = 118816.445
$ws.Range("K34").Value = 9379413
$ws.Range("L34").Value = 356449.335
$ws.Range("M34").Value = -9379329
$ws.Range("N34").Value = -356617.335

$ws.Range("H44").Value = 1876.2778
$ws.Range("J44").Value = 4406.6665
$ws.Range("L44").Value = 13219.9995
$ws.Range("N44").Value = -14015.9995

$ws.Range("H46").Value = 2549.8333
$ws.Range("I46").Value = 474.75
$ws.Range("K46").Value = 1424.25
$ws.Range("M46").Value = -1333.25

$ws.Range("H68").Value = 4262.647
$ws.Range("I68").Value = 1919.5
$ws.Range("K68").Value = 5758.5
$ws.Range("M68").Value = -4947.5

$ws.Range("H71").Value = 4262.647
$ws.Range("I71").Value = 1919.5
$ws.Range("K71").Value = 17275.5
$ws.Range("M71").Value = -13219.5

$ws.Range("H113").Value = 1254.7188
$ws.Range("J113").Value = 1555.15
$ws.Range("L113").Value = 4665.450000000001
$ws.Range("N113").Value = -9005.450000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5336.1763
$ws.Range("I70").Value = 6127
$ws.Range("K70").Value = 6127
$ws.Range("M70").Value = -5857

$ws.Range("H73").Value = 5336.1763
$ws.Range("I73").Value = 6127
$ws.Range("K73").Value = 6127
$ws.Range("M73").Value = -5191

$ws.Range("H80").Value = 12957.3
$ws.Range("I80").Value = 25944.75
$ws.Range("J80").Value = 4299
$ws.Range("K80").Value = 25944.75
$ws.Range("L80").Value = 4299
$ws.Range("M80").Value = -24946.75
$ws.Range("N80").Value = -6295

$ws.Range("H83").Value = 12957.3
$ws.Range("I83").Value = 25944.75
$ws.Range("J83").Value = 4299
$ws.Range("K83").Value = 129723.75
$ws.Range("L83").Value = 21495
$ws.Range("M83").Value = -124731.75
$ws.Range("N83").Value = -31479

$ws.Range("H97").Value = 8512.414000000001
$ws.Range("I97").Value = 10235.228
$ws.Range("K97").Value = 10235.228
$ws.Range("M97").Value = -9739.227999999999

$ws.Range("H132").Value = 4241.4443
$ws.Range("I132").Value = 3373.4707
$ws.Range("K132").Value = 10120.4121
$ws.Range("M132").Value = -7590.4121

$ws.Range("H133").Value = 39999
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 94591.44500000001
$ws.Range("J136").Value = 102665.375
$ws.Range("L136").Value = 307996.125
$ws.Range("N136").Value = -313096.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10728
$ws.Range("J22").Value = 6924.7
$ws.Range("L22").Value = 6924.7
$ws.Range("N22").Value = -7514.7

$ws.Range("H27").Value = 10728
$ws.Range("J27").Value = 6924.7
$ws.Range("L27").Value = 6924.7
$ws.Range("N27").Value = -7138.7

$ws.Range("H40").Value = 40219.43
$ws.Range("I40").Value = 45913.688
$ws.Range("K40").Value = 45913.688
$ws.Range("M40").Value = -45777.688

$ws.Range("H98").Value = 100355
$ws.Range("J98").Value = 100355
$ws.Range("L98").Value = 100355
$ws.Range("N98").Value = -106345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 32293.588
$ws.Range("I100").Value = 16539.4
$ws.Range("J100").Value = 150450
$ws.Range("K100").Value = 33078.8
$ws.Range("L100").Value = 300900
$ws.Range("M100").Value = -32537.8
$ws.Range("N100").Value = -301982

$ws.Range("H122").Value = 7017.0625
$ws.Range("I122").Value = 2325.7144
$ws.Range("J122").Value = 10665.889
$ws.Range("K122").Value = 6977.1432
$ws.Range("L122").Value = 31997.667
$ws.Range("M122").Value = -4527.1432
$ws.Range("N122").Value = -36897.667

$ws.Range("H132").Value = 24125.104
$ws.Range("I132").Value = 33363.445
$ws.Range("K132").Value = 100090.335
$ws.Range("M132").Value = -97560.33499999999

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
